# Scheduled-runner refresh of currentAveragePrice / Leve price / profit columns
# (H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Generated from the authoritative OOXML diff; values are written as plain
# numeric literals (source cells hold no formulas). Where the refreshed
# LevePriceNQ (K) / LevePriceHQ (L) becomes 0, the corresponding profit cell
# (M / N) is removed entirely (ClearContents), matching the diff; likewise new
# profit cells are created where the diff adds them.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1127
$ws.Range("I11").Value = 1127
$ws.Range("K11").Value = 1127
$ws.Range("M11").Value = -987
$ws.Range("H33").Value = 268.57895
$ws.Range("I33").Value = 157.16667
$ws.Range("K33").Value = 157.16667
$ws.Range("M33").Value = 71.83332999999999
$ws.Range("H39").Value = 5845.909
$ws.Range("I39").Value = 2151.5
$ws.Range("J39").Value = 7957
$ws.Range("K39").Value = 6454.5
$ws.Range("L39").Value = 23871
$ws.Range("M39").Value = -6158.5
$ws.Range("N39").Value = -24463
$ws.Range("H42").Value = 1960.4166
$ws.Range("I42").Value = 1331.7778
$ws.Range("K42").Value = 3995.3334
$ws.Range("M42").Value = -3765.3334
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H76").Value = 4159.1665
$ws.Range("I76").Value = 3351.6667
$ws.Range("K76").Value = 3351.6667
$ws.Range("M76").Value = -3036.6667
$ws.Range("H79").Value = 4159.1665
$ws.Range("I79").Value = 3351.6667
$ws.Range("K79").Value = 3351.6667
$ws.Range("M79").Value = -2259.6667
$ws.Range("H92").Value = 789.9048
$ws.Range("I92").Value = 749.2222
$ws.Range("K92").Value = 749.2222
$ws.Range("M92").Value = 498.7778
$ws.Range("H99").Value = 2679.5454
$ws.Range("I99").Value = 247.125
$ws.Range("J99").Value = 9166
$ws.Range("K99").Value = 741.375
$ws.Range("L99").Value = 27498
$ws.Range("M99").Value = 756.625
$ws.Range("N99").Value = -30494
$ws.Range("H132").Value = 3174.9524
$ws.Range("I132").Value = 3142.647
$ws.Range("J132").Value = 3312.25
$ws.Range("K132").Value = 9427.940999999999
$ws.Range("L132").Value = 9936.75
$ws.Range("M132").Value = -6897.940999999999
$ws.Range("N132").Value = -14996.75
$ws.Range("H135").Value = 150000770
$ws.Range("I135").Value = 62500652
$ws.Range("K135").Value = 562505868
$ws.Range("M135").Value = -562503333
$ws.Range("H137").Value = 3218.5334
$ws.Range("I137").Value = 2490.2856
$ws.Range("K137").Value = 7470.8568
$ws.Range("M137").Value = -4920.8568
$ws.Range("H138").Value = 1750.4667
$ws.Range("I138").Value = 1290.3793
$ws.Range("K138").Value = 3871.1379
$ws.Range("M138").Value = 1268.8621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2261.2632
$ws.Range("I32").Value = 2241.2974
$ws.Range("K32").Value = 2241.2974
$ws.Range("M32").Value = -1954.2974
$ws.Range("H34").Value = 28755.5
$ws.Range("I34").Value = 17512
$ws.Range("K34").Value = 17512
$ws.Range("M34").Value = -17241
$ws.Range("H45").Value = 7283.4375
$ws.Range("I45").Value = 7830.4614
$ws.Range("K45").Value = 7830.4614
$ws.Range("M45").Value = -7453.4614
$ws.Range("H61").Value = 62501904
$ws.Range("I61").Value = 71430430
$ws.Range("K61").Value = 71430430
$ws.Range("M61").Value = -71430218
$ws.Range("H74").Value = 33336000
$ws.Range("I74").Value = 38464440
$ws.Range("K74").Value = 38464440
$ws.Range("M74").Value = -38463566
$ws.Range("H77").Value = 33336000
$ws.Range("I77").Value = 38464440
$ws.Range("K77").Value = 192322200
$ws.Range("M77").Value = -192317832
$ws.Range("H122").Value = 5740.52
$ws.Range("I122").Value = 4125.7
$ws.Range("J122").Value = 12199.8
$ws.Range("K122").Value = 12377.1
$ws.Range("L122").Value = 36599.39999999999
$ws.Range("M122").Value = -9927.099999999999
$ws.Range("N122").Value = -41499.39999999999
$ws.Range("H132").Value = 3451143
$ws.Range("I132").Value = 3451143
$ws.Range("K132").Value = 10353429
$ws.Range("M132").Value = -10350899
$ws.Range("H136").Value = 62501904
$ws.Range("I136").Value = 71430430
$ws.Range("K136").Value = 214291290
$ws.Range("M136").Value = -214288740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4106.3477
$ws.Range("I86").Value = 4596.7334
$ws.Range("J86").Value = 3186.875
$ws.Range("K86").Value = 4596.7334
$ws.Range("L86").Value = 3186.875
$ws.Range("M86").Value = -3473.7334
$ws.Range("N86").Value = -5432.875
$ws.Range("H89").Value = 4106.3477
$ws.Range("I89").Value = 4596.7334
$ws.Range("J89").Value = 3186.875
$ws.Range("K89").Value = 22983.667
$ws.Range("L89").Value = 15934.375
$ws.Range("M89").Value = -17367.667
$ws.Range("N89").Value = -27166.375
$ws.Range("H122").Value = 79999
$ws.Range("J122").Value = 79999
$ws.Range("L122").Value = 79999
$ws.Range("N122").Value = -89799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 27749.5
$ws.Range("I22").Value = 50249
$ws.Range("K22").Value = 50249
$ws.Range("M22").Value = -49899
$ws.Range("H31").Value = 6673.0835
$ws.Range("I31").Value = 3907.7
$ws.Range("J31").Value = 20500
$ws.Range("K31").Value = 3907.7
$ws.Range("L31").Value = 20500
$ws.Range("M31").Value = -3612.7
$ws.Range("N31").Value = -21090
$ws.Range("H34").Value = 6673.0835
$ws.Range("I34").Value = 3907.7
$ws.Range("J34").Value = 20500
$ws.Range("K34").Value = 3907.7
$ws.Range("L34").Value = 20500
$ws.Range("M34").Value = -3705.7
$ws.Range("N34").Value = -20904
$ws.Range("H106").Value = 15000
$ws.Range("J106").Value = 15000
$ws.Range("L106").Value = 15000
$ws.Range("N106").Value = -17524
$ws.Range("H107").Value = 34128.8
$ws.Range("I107").Value = 618.7917
$ws.Range("J107").Value = 168168.83
$ws.Range("K107").Value = 618.7917
$ws.Range("L107").Value = 168168.83
$ws.Range("M107").Value = 1301.2083
$ws.Range("N107").Value = -172008.83
$ws.Range("H132").Value = 333334900
$ws.Range("I132").Value = 333334900
$ws.Range("K132").Value = 1000004700
$ws.Range("M132").Value = -1000002170
$ws.Range("H141").Value = 83763.336
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 93410
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 93410
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -103770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 184
$ws.Range("I51").Value = 184
$ws.Range("K51").Value = 552
$ws.Range("M51").Value = -92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3473.6316
$ws.Range("I102").Value = 3473.6316
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3473.6316
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1851.6316
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 5200.222
$ws.Range("I122").Value = 2199.1667
$ws.Range("K122").Value = 6597.500100000001
$ws.Range("M122").Value = -4147.500100000001
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H134").Value = 58464.2
$ws.Range("J134").Value = 58464.2
$ws.Range("L134").Value = 175392.6
$ws.Range("N134").Value = -180462.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4500
$ws.Range("I7").Value = 4500
$ws.Range("K7").Value = 4500
$ws.Range("M7").Value = -4388
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H60").Value = 18030.5
$ws.Range("I60").Value = 16000
$ws.Range("J60").Value = 20061
$ws.Range("K60").Value = 16000
$ws.Range("L60").Value = 20061
$ws.Range("M60").Value = -15491
$ws.Range("N60").Value = -21079
$ws.Range("H98").Value = 35177.5
$ws.Range("I98").Value = 30000
$ws.Range("J98").Value = 40355
$ws.Range("K98").Value = 30000
$ws.Range("L98").Value = 40355
$ws.Range("M98").Value = -27005
$ws.Range("N98").Value = -46345
$ws.Range("H100").Value = 16634487
$ws.Range("I100").Value = 19960904
$ws.Range("K100").Value = 19960904
$ws.Range("M100").Value = -19960363
$ws.Range("H122").Value = 5847.2354
$ws.Range("I122").Value = 3707.4285
$ws.Range("J122").Value = 15833
$ws.Range("K122").Value = 11122.2855
$ws.Range("L122").Value = 47499
$ws.Range("M122").Value = -8672.2855
$ws.Range("N122").Value = -52399
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030
$ws.Range("H132").Value = 10425361
$ws.Range("I132").Value = 13166456
$ws.Range("K132").Value = 39499368
$ws.Range("M132").Value = -39496838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 500000
$ws.Range("I5").Value = 500000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -499888
$ws.Range("N5").ClearContents()
$ws.Range("H40").Value = 23329.666
$ws.Range("I40").Value = 19995
$ws.Range("J40").Value = 29999
$ws.Range("K40").Value = 19995
$ws.Range("L40").Value = 29999
$ws.Range("M40").Value = -19846
$ws.Range("N40").Value = -30297
$ws.Range("H126").Value = 3030.3215
$ws.Range("I126").Value = 3127.6
$ws.Range("K126").Value = 9382.799999999999
$ws.Range("M126").Value = -6912.799999999999
$ws.Range("H132").Value = 12507111
$ws.Range("I132").Value = 17860608
$ws.Range("K132").Value = 53581824
$ws.Range("M132").Value = -53579294
